# Edit: "Edited language to get as close to 9th grade level as possible."
#
# NOTE on ordering: in this COM-interop runtime, Range.Find.Execute with
# Replace=wdReplaceAll (2) walks forward from the range's Start position to
# the end of the *whole document* story (it does not stop at the range's own
# End). So body-text edits that live further down the document are done
# first, and the Title (which sits at document position 0 and also contains
# matching text) is edited last, using Replace=wdReplaceOne (1) so it only
# touches its own single occurrence.
#
# Also: touching $d.Tables / Cell ranges before finishing the
# $d.Paragraphs.Item(n) based edits desyncs later paragraph-index lookups in
# this runtime, so the Title edit below is done purely through
# $d.Content.Find (document Range), never via $d.Tables.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Congratulations" paragraph (body): simplify "File" -> "complete" and
#    tidy up the surrounding punctuation/spacing.
# ---------------------------------------------------------------------

# "...forms you need to File an affidavit of address verification." ->
# "...forms you need to complete an affidavit of address verification."
$p4 = $d.Paragraphs.Item(4).Range
$p4.Find.Execute("File an affidavit of address verification", $true, $false, $false, $false, $false, $true, 0, $false, "complete an affidavit of address verification", 2) | Out-Null

# Collapse the double space before the "{% if ... %}" tag down to one space.
$p4b = $d.Paragraphs.Item(4).Range
$p4b.Find.Execute("Verification  {% if", $true, $false, $false, $false, $false, $true, 0, $false, "Verification {% if", 2) | Out-Null

# " v {{ " -> " v. {{ "  (add the period after "v")
$p4c = $d.Paragraphs.Item(4).Range
$p4c.Find.Execute(" v {{", $true, $false, $false, $false, $false, $true, 0, $false, " v. {{", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Heading2 "To file your affidavit right away" -> "File your affidavit
#    right away"
# ---------------------------------------------------------------------
$d.Content.Find.Execute("To file your ", $true, $false, $false, $false, $false, $true, 0, $false, "File your ", 2) | Out-Null

# ---------------------------------------------------------------------
# 3) Title (table banner): "File an affidavit of address verification" ->
#    "Complete an affidavit of address verification"
#    (Replace=wdReplaceOne so only the first/title occurrence is touched,
#    now that the body copy above no longer contains this phrase.)
# ---------------------------------------------------------------------
$d.Content.Find.Execute("File an affidavit of address verification", $true, $false, $false, $false, $false, $true, 0, $false, "Complete an affidavit of address verification", 1) | Out-Null

# ---------------------------------------------------------------------
# 4) Remove the stray empty paragraph right after the "Congratulations"
#    paragraph (between it and the section-break paragraph).
# ---------------------------------------------------------------------
$d.Paragraphs.Item(5).Range.Delete()
